$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.147.68"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.25%  "

$ws.Range("D3").Value = "'3.403.97"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.61%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'581.88"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.46%  "

$ws.Range("D6").Value = "'178.24"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.71%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("E8").Value = "  +0.43%  "

$ws.Range("E9").Value = "  +7.91%  "

$ws.Range("D10").Value = "'0.585"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("D11").Value = "'48.28"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.77%  "

$ws.Range("E12").Value = "  +3.05%  "

$ws.Range("D13").Value = "'678.71"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("D14").Value = "'3.953.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.40%  "

$ws.Range("E15").Value = "  +1.86%  "

$ws.Range("D16").Value = "'69.261.38"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.36%  "

$ws.Range("D17").Value = "'3.403.30"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.07%  "

$ws.Range("E18").Value = "  +0.49%  "

$ws.Range("E19").Value = "  +1.47%  "

$ws.Range("E20").Value = "  +0.96%  "

$ws.Range("D21").Value = "'0.910"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.83%  "

$ws.Range("D22").Value = "'5.37"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.71%  "

$ws.Range("D23").Value = "'17.02"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("D24").Value = "'100.75"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.60%  "

$ws.Range("D25").Value = "'3.89"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.29%  "

$ws.Range("E26").Value = "  -0.30%  "

$ws.Range("D27").Value = "'9.66"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.03%  "

$ws.Range("D28").Value = "'33.48"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.55%  "

$ws.Range("E29").Value = "  +2.56%  "

$ws.Range("D30").Value = "'6.85"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("E31").Value = "  +10.54%  "

$ws.Range("D32").Value = "'555.15"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.58%  "

$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("D35").Value = "'57.99"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("E36").Value = "  +0.10%  "

$ws.Range("D37").Value = "'3.609.60"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.94%  "

$ws.Range("E38").Value = "  +0.62%  "

$ws.Range("D39").Value = "'34.96"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.06%  "

$ws.Range("D40").Value = "'0.0₃0738"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +10.14%  "

$ws.Range("E41").Value = "  +3.49%  "

$ws.Range("E42").Value = "  +2.86%  "

$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0424"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.31%  "

$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").Value = "'0.334"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.08%  "

$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'2.66"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.65%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.129"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'1.40"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.72%  "

$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "'1.00"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.19%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'130.99"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.79%  "

$ws.Range("B50").Value = "CoreDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D50").Value = "'2.68"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.44%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").Value = "'7.38"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.22%  "
